$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.771.35'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '2.347.32'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '545.28'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.85%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.84'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.573'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +7.65%  '
$ws.Range('D9').Value = '2.344.63'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.102'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.40'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.17%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.153'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.359'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +6.38%  '
$ws.Range('D14').Value = '2.761.12'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.59'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '57.736.46'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').Value = '2.344.92'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.63'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.12%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '334.89'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.71'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '62.08'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.55'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('E29').Value = '  +6.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.77'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '170.47'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('D32').Value = '0.0₃0733'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.15'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +16.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.49'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.17'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.05%  '
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.62'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.50%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '39.09'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '148.11'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.62'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '284.87'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0942'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.17'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +6.16%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0506'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.96%  '
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.51'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.94%  '
